$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab date from 2025-10-14 to 2025-10-15
$ws.Name = "PickAndPlace_PCB_2025-10-15"

# Row 15 previously held the LED1 part; it now holds the R1 resistor data,
# with an updated rotation (0) and comment (150Ω).
$ws.Range("A15").Value = "R1"
$ws.Range("B15").Value = "RS-06K151JT_C140047"
$ws.Range("C15").Value = "R1206"
$ws.Range("D15").Value = "74.93mm"
$ws.Range("E15").Value = "-38.227mm"
$ws.Range("F15").Value = "74.93mm"
$ws.Range("G15").Value = "-38.227mm"
$ws.Range("H15").Value = "73.451mm"
$ws.Range("I15").Value = "-38.227mm"
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = "150Ω"

# Row 16 previously held the R1 part; it now holds the LED data with
# refreshed device/footprint/pad-X values and comment.
$ws.Range("A16").Value = "LED"
$ws.Range("B16").Value = "YLED1206R"
$ws.Range("C16").Value = "LED1206-FD"
$ws.Range("D16").Value = "80.264mm"
$ws.Range("E16").Value = "-38.227mm"
$ws.Range("F16").Value = "80.264mm"
$ws.Range("G16").Value = "-38.227mm"
$ws.Range("H16").Value = "78.747mm"
$ws.Range("I16").Value = "-38.227mm"
$ws.Range("N16").Value = "YLED1206R"
